$wb = $excel.ActiveWorkbook

# --- Sheet "Stato Attuale" (current fleet status: targa | operatore | data_assegnazione) ---
$ws1 = $wb.Worksheets.Item("Stato Attuale")

# Targa GL777AD (row 19): reassigned from GAIA.MARTI to "DA ASSEGNARE(MANDARE IN ASSISTENZA)"
# on 2026-02-03 (previously had no assignment date recorded).
$ws1.Range("B19").Value = "DA ASSEGNARE(MANDARE IN ASSISTENZA)"
$ws1.Range("C19").NumberFormat = "@"
$ws1.Range("C19").Value = "2026-02-03"
$ws1.Range("C19").Style = "Normal"

# Targa GL594TH (row 76): reassigned from "DA ASSEGNARE" to
# "ROCCO ZACCAGNIGNO (MOMENTANEO, SUA IN ASSISTENZA)"; assignment date moves
# from 2026-02-02 to 2026-02-03.
$ws1.Range("B76").Value = "ROCCO ZACCAGNIGNO (MOMENTANEO, SUA IN ASSISTENZA)"
$ws1.Range("C76").NumberFormat = "@"
$ws1.Range("C76").Value = "2026-02-03"
$ws1.Range("C76").Style = "Normal"

# --- Sheet "Storico Passaggi" (change history: Targa | Operatore_Precedente | Nuovo_Operatore | Data_Cambio) ---
$ws2 = $wb.Worksheets.Item("Storico Passaggi")

# Row 2 now logs the GL777AD hand-off.
$ws2.Range("A2").Value = "GL777AD"
$ws2.Range("B2").Value = "GAIA.MARTI"
$ws2.Range("C2").Value = "DA ASSEGNARE(MANDARE IN ASSISTENZA)"
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "2026-02-03"
$ws2.Range("D2").Style = "Normal"

# Row 3 is new and logs the GL594TH hand-off.
$ws2.Range("A3").Value = "GL594TH"
$ws2.Range("B3").Value = "DA ASSEGNARE"
$ws2.Range("C3").Value = "ROCCO ZACCAGNIGNO (MOMENTANEO, SUA IN ASSISTENZA)"
$ws2.Range("D3").NumberFormat = "@"
$ws2.Range("D3").Value = "2026-02-03"
$ws2.Range("D3").Style = "Normal"
